# Applies scheduled market-price / leve-profit data refresh to the Fenrir_Profits workbook.
# Each sheet (ALC, ARM, BSM, CRP, GSM, LTW, WVR) gets updated currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H:N) for the rows whose source data changed.
# CUL is untouched by this run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 114.44444
$ws.Range("I55").Value = 106
$ws.Range("J55").Value = 125
$ws.Range("K55").Value = 106
$ws.Range("L55").Value = 125
$ws.Range("M55").Value = 108
$ws.Range("N55").Value = -553
# Row 74
$ws.Range("H74").Value = 5877.143
$ws.Range("I74").Value = 6912.857
$ws.Range("J74").Value = 4841.4287
$ws.Range("K74").Value = 6912.857
$ws.Range("L74").Value = 4841.4287
$ws.Range("M74").Value = -5976.857
$ws.Range("N74").Value = -6713.4287
# Row 76
$ws.Range("H76").Value = 4360.16
$ws.Range("I76").Value = 4028.5715
$ws.Range("J76").Value = 4489.1113
$ws.Range("K76").Value = 4028.5715
$ws.Range("L76").Value = 4489.1113
$ws.Range("M76").Value = -3713.5715
$ws.Range("N76").Value = -5119.1113
# Row 77
$ws.Range("H77").Value = 5877.143
$ws.Range("I77").Value = 6912.857
$ws.Range("J77").Value = 4841.4287
$ws.Range("K77").Value = 34564.285
$ws.Range("L77").Value = 24207.1435
$ws.Range("M77").Value = -29884.285
$ws.Range("N77").Value = -33567.14350000001
# Row 79
$ws.Range("H79").Value = 4360.16
$ws.Range("I79").Value = 4028.5715
$ws.Range("J79").Value = 4489.1113
$ws.Range("K79").Value = 4028.5715
$ws.Range("L79").Value = 4489.1113
$ws.Range("M79").Value = -2936.5715
$ws.Range("N79").Value = -6673.1113
# Row 116
$ws.Range("H116").Value = 299448.03
$ws.Range("I116").Value = 4220.2
$ws.Range("K116").Value = 4220.2
$ws.Range("M116").Value = -778.1999999999998
# Row 135
$ws.Range("H135").Value = 1340.5143
$ws.Range("I135").Value = 1138.1
$ws.Range("K135").Value = 10242.9
$ws.Range("M135").Value = -7707.9
# Row 137
$ws.Range("H137").Value = 1180.1471
$ws.Range("I137").Value = 909.86365
$ws.Range("J137").Value = 1675.6666
$ws.Range("K137").Value = 2729.59095
$ws.Range("L137").Value = 5026.9998
$ws.Range("M137").Value = -179.5909499999998
$ws.Range("N137").Value = -10126.9998
# Row 138
$ws.Range("H138").Value = 1471
$ws.Range("I138").Value = 947.04346
$ws.Range("J138").Value = 2475.25
$ws.Range("K138").Value = 2841.13038
$ws.Range("L138").Value = 7425.75
$ws.Range("M138").Value = 2298.86962
$ws.Range("N138").Value = -17705.75
# Row 141
$ws.Range("H141").Value = 9180.294
$ws.Range("I141").Value = 9096.666999999999
$ws.Range("J141").Value = 9381
$ws.Range("K141").Value = 27290.001
$ws.Range("L141").Value = 28143
$ws.Range("M141").Value = -22110.001
$ws.Range("N141").Value = -38503

# ---------------------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1198.4166
$ws.Range("I45").Value = 1221.75
$ws.Range("J45").Value = 1151.75
$ws.Range("K45").Value = 1221.75
$ws.Range("L45").Value = 1151.75
$ws.Range("M45").Value = -844.75
$ws.Range("N45").Value = -1905.75
# Row 61
$ws.Range("H61").Value = 3823.8718
$ws.Range("I61").Value = 4448.72
$ws.Range("J61").Value = 2708.0715
$ws.Range("K61").Value = 4448.72
$ws.Range("L61").Value = 2708.0715
$ws.Range("M61").Value = -4236.72
$ws.Range("N61").Value = -3132.0715
# Row 63
$ws.Range("H63").Value = 1854683
$ws.Range("I63").Value = 4168202.8
$ws.Range("J63").Value = 3867.2
$ws.Range("K63").Value = 4168202.8
$ws.Range("L63").Value = 3867.2
$ws.Range("M63").Value = -4167516.8
$ws.Range("N63").Value = -5239.2
# Row 66
$ws.Range("H66").Value = 1854683
$ws.Range("I66").Value = 4168202.8
$ws.Range("J66").Value = 3867.2
$ws.Range("K66").Value = 20841014
$ws.Range("L66").Value = 19336
$ws.Range("M66").Value = -20837582
$ws.Range("N66").Value = -26200
# Row 74
$ws.Range("H74").Value = 1337.6471
$ws.Range("I74").Value = 653
$ws.Range("J74").Value = 1548.3077
$ws.Range("K74").Value = 653
$ws.Range("L74").Value = 1548.3077
$ws.Range("M74").Value = 221
$ws.Range("N74").Value = -3296.3077
# Row 77
$ws.Range("H77").Value = 1337.6471
$ws.Range("I77").Value = 653
$ws.Range("J77").Value = 1548.3077
$ws.Range("K77").Value = 3265
$ws.Range("L77").Value = 7741.538500000001
$ws.Range("M77").Value = 1103
$ws.Range("N77").Value = -16477.5385
# Row 132
$ws.Range("H132").Value = 2360389.8
$ws.Range("I132").Value = 4808934.5
$ws.Range("J132").Value = 2532.1853
$ws.Range("K132").Value = 14426803.5
$ws.Range("L132").Value = 7596.5559
$ws.Range("M132").Value = -14424273.5
$ws.Range("N132").Value = -12656.5559
# Row 134
$ws.Range("H134").Value = 25750
$ws.Range("J134").Value = 25750
$ws.Range("L134").Value = 25750
$ws.Range("N134").Value = -35890
# Row 136
$ws.Range("H136").Value = 3823.8718
$ws.Range("I136").Value = 4448.72
$ws.Range("J136").Value = 2708.0715
$ws.Range("K136").Value = 13346.16
$ws.Range("L136").Value = 8124.2145
$ws.Range("M136").Value = -10796.16
$ws.Range("N136").Value = -13224.2145

# ---------------------------------------------------------------------------
# Sheet: BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1175.5
$ws.Range("I22").Value = 1175.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1175.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1002.5
$ws.Range("N22").ClearContents()
# Row 105
$ws.Range("H105").Value = 1968.1428
$ws.Range("I105").Value = 1701
$ws.Range("J105").Value = 2168.5
$ws.Range("K105").Value = 1701
$ws.Range("L105").Value = 2168.5
$ws.Range("M105").Value = 46
$ws.Range("N105").Value = -5662.5
# Row 134
$ws.Range("H134").Value = 9539594
$ws.Range("I134").Value = 13354587
$ws.Range("J134").Value = 2110.9
$ws.Range("K134").Value = 40063761
$ws.Range("L134").Value = 6332.700000000001
$ws.Range("M134").Value = -40061226
$ws.Range("N134").Value = -11402.7

# ---------------------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3692681.8
$ws.Range("I58").Value = 5995382
$ws.Range("J58").Value = 8361.532999999999
$ws.Range("K58").Value = 5995382
$ws.Range("L58").Value = 8361.532999999999
$ws.Range("M58").Value = -5995179
$ws.Range("N58").Value = -8767.532999999999
# Row 62
$ws.Range("H62").Value = 5600.4
$ws.Range("I62").Value = 7860.8
$ws.Range("J62").Value = 3340
$ws.Range("K62").Value = 7860.8
$ws.Range("L62").Value = 3340
$ws.Range("M62").Value = -7236.8
$ws.Range("N62").Value = -4588
# Row 65
$ws.Range("H65").Value = 5600.4
$ws.Range("I65").Value = 7860.8
$ws.Range("J65").Value = 3340
$ws.Range("K65").Value = 39304
$ws.Range("L65").Value = 16700
$ws.Range("M65").Value = -36184
$ws.Range("N65").Value = -22940
# Row 99
$ws.Range("H99").Value = 4971.3335
$ws.Range("I99").Value = 7703
$ws.Range("J99").Value = 2786
$ws.Range("K99").Value = 7703
$ws.Range("L99").Value = 2786
$ws.Range("M99").Value = -6205
$ws.Range("N99").Value = -5782
# Row 126
$ws.Range("H126").Value = 4971.3335
$ws.Range("I126").Value = 7703
$ws.Range("J126").Value = 2786
$ws.Range("K126").Value = 23109
$ws.Range("L126").Value = 8358
$ws.Range("M126").Value = -20639
$ws.Range("N126").Value = -13298
# Row 134
$ws.Range("H134").Value = 8682429
$ws.Range("I134").Value = 14707771
$ws.Range("J134").Value = 3291333.2
$ws.Range("K134").Value = 44123313
$ws.Range("L134").Value = 9873999.600000001
$ws.Range("M134").Value = -44120778
$ws.Range("N134").Value = -9879069.600000001
# Row 136
$ws.Range("H136").Value = 3692681.8
$ws.Range("I136").Value = 5995382
$ws.Range("J136").Value = 8361.532999999999
$ws.Range("K136").Value = 17986146
$ws.Range("L136").Value = 25084.599
$ws.Range("M136").Value = -17983596
$ws.Range("N136").Value = -30184.599

# ---------------------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# Row 70
$ws.Range("H70").Value = 203981.6
$ws.Range("I70").Value = 502454
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 502454
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -502184
$ws.Range("N70").Value = -5540
# Row 73
$ws.Range("H73").Value = 203981.6
$ws.Range("I73").Value = 502454
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 502454
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -501518
$ws.Range("N73").Value = -6872
# Row 80
$ws.Range("H80").Value = 3200
$ws.Range("I80").Value = 2800
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 2800
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -1802
$ws.Range("N80").Value = -5996
# Row 83
$ws.Range("H83").Value = 3200
$ws.Range("I83").Value = 2800
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 14000
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -9008
$ws.Range("N83").Value = -29984
# Row 132
$ws.Range("H132").Value = 37079360
$ws.Range("I132").Value = 83418880
$ws.Range("J132").Value = 7745.4
$ws.Range("K132").Value = 250256640
$ws.Range("L132").Value = 23236.2
$ws.Range("M132").Value = -250254110
$ws.Range("N132").Value = -28296.2

# ---------------------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1971.4286
$ws.Range("I7").Value = 1971.4286
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1971.4286
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1859.4286
$ws.Range("N7").ClearContents()
# Row 126
$ws.Range("H126").Value = 1971.4286
$ws.Range("I126").Value = 1971.4286
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5914.2858
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3444.2858
$ws.Range("N126").ClearContents()
# Row 136
$ws.Range("H136").Value = 3860.35
$ws.Range("I136").Value = 4093
$ws.Range("J136").Value = 3317.5
$ws.Range("K136").Value = 12279
$ws.Range("L136").Value = 9952.5
$ws.Range("M136").Value = -9729
$ws.Range("N136").Value = -15052.5

# ---------------------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 24250
$ws.Range("J64").Value = 24250
$ws.Range("L64").Value = 24250
$ws.Range("N64").Value = -24746
# Row 67
$ws.Range("H67").Value = 24250
$ws.Range("J67").Value = 24250
$ws.Range("L67").Value = 24250
$ws.Range("N67").Value = -25966
# Row 70
$ws.Range("H70").Value = 17491.25
$ws.Range("J70").Value = 17491.25
$ws.Range("L70").Value = 17491.25
$ws.Range("N70").Value = -18121.25
# Row 73
$ws.Range("H73").Value = 17491.25
$ws.Range("J73").Value = 17491.25
$ws.Range("L73").Value = 17491.25
$ws.Range("N73").Value = -19675.25
# Row 122
$ws.Range("H122").Value = 955.7857
$ws.Range("I122").Value = 777.1
$ws.Range("J122").Value = 1402.5
$ws.Range("K122").Value = 2331.3
$ws.Range("L122").Value = 4207.5
$ws.Range("M122").Value = 118.6999999999998
$ws.Range("N122").Value = -9107.5

